$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '29.941.54'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.633.71'
$ws.Range("E3").Value = '  +1.82%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '214.72'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '28.74'
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").Value = '0.0903'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '1.867.48'
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("D13").Value = '1.634.26'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '9.28'
$ws.Range("E15").Value = '  +13.77%  '
$ws.Range("D16").Value = '29.962.42'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '3.84'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '64.16'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '241.03'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '0.0₃0702'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +2.36%  '
$ws.Range("D23").Value = '9.80'
$ws.Range("E23").Value = '  +3.28%  '
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '157.72'
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").Value = '15.48'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").Value = '6.56'
$ws.Range("E28").Value = '  +1.13%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("E31").Value = '  +3.58%  '
$ws.Range("D32").Value = '3.39'
$ws.Range("E32").Value = '  +4.42%  '
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").Value = '1.431.62'
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").Value = '  +5.12%  '
$ws.Range("E36").Value = '  -1.76%  '
$ws.Range("D37").Value = '2.77'
$ws.Range("E37").Value = '  -2.42%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").Value = '75.66'
$ws.Range("E40").Value = '  +11.76%  '
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").Value = '1.99'
$ws.Range("E42").Value = '  +1.37%  '
$ws.Range("D43").Value = '0.829'
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("D44").Value = '0.0496'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").Value = '51.35'
$ws.Range("E47").Value = '  -7.15%  '
$ws.Range("D48").Value = '5.35'
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("D49").Value = '1.774.52'
$ws.Range("E49").Value = '  +1.97%  '
$ws.Range("E50").Value = '  +9.93%  '
$ws.Range("D51").Value = '90.41'
$ws.Range("E51").Value = '  +4.18%  '